$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# The Java Script function in G2 was updated: FullScreen is now enabled (value=1)
# and a new isWindowsKey=0 config step was added.
$ws.Range("G2").Value2 = "wait(3);`nPullConfigxml;`nChangeConfigxml(Configuration/Applications/Application/General,StartPage,<StartPage value=`"http://127.0.0.1:8082/app/`" name=`"Menu`"/>);`nChangeConfigxml(Configuration,WebServer,<WebServer>endl  <Enabled VALUE=`"1`"/>endl  <Port VALUE=`"8082`"/>endl  <WebFolder VALUE=`"\\auto\\RE_2.2\`"/>endl  <Public VALUE=`"1`"/>endl</WebServer>endl);`nChangeConfigxml(Configuration/Screen,FullScreen,<FullScreen value=`"1`"/>);`nChangeConfigxml(Configuration/Applications/Application,isWindowsKey,<isWindowsKey value=`"0`"/>);`nChangeConfigxml(Configuration/DeviceKeys,FunctionKeysCapturable,<FunctionKeysCapturable value=`"1`"/>);`nPushConfigxml;"

# The old per-row "Pass" markers in column J were cleared out.
$ws.Range("J2:J31").ClearContents()

# Row 2 grew taller to fit the extra line of script text.
$ws.Rows.Item(2).RowHeight = 332.25
